$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("A2").Value = 111813975
$ws.Range("B2").Value = 89423
$ws.Range("D2").Value = 'NT'
$ws.Range("E2").Value = 5432
$ws.Range("F2").Value = 'Granticka'
$ws.Range("G2").Value = 'Porodaedalea chrysoloma'
$ws.Range("H2").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q2").Value = 540643.7191088985
$ws.Range("R2").Value = 7247516.737328541

# ---- Row 3 ----
$ws.Range("A3").Value = 111814119
$ws.Range("Q3").Value = 540683.0369185829
$ws.Range("R3").Value = 7247576.171207689

# ---- Row 4 ----
$ws.Range("A4").Value = 111813707
$ws.Range("Q4").Value = 540647.037727406
$ws.Range("R4").Value = 7247579.013394679
$ws.Range("AJ4").Value = 'gran'
$ws.Range("AK4").Value = 'Picea abies'
$ws.Range("AO4").Value = 'Picea abies'

# ---- Row 5 ----
$ws.Range("A5").Value = 111813872
$ws.Range("B5").Value = 56398
$ws.Range("D5").Value = 'NT'
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = 'Tretåig hackspett'
$ws.Range("G5").Value = 'Picoides tridactylus'
$ws.Range("H5").Value = '(Linnaeus, 1758)'
$ws.Range("J5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = 'färska spår'
$ws.Range("Q5").Value = 540557.5018987871
$ws.Range("R5").Value = 7247552.715308581
$ws.Range("AF5").Value = ""

# ---- Row 7 ----
$ws.Range("A7").Value = 111814135
$ws.Range("B7").Value = 90087
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 3298
$ws.Range("F7").Value = 'Trådticka'
$ws.Range("G7").Value = 'Climacocystis borealis'
$ws.Range("H7").Value = '(Fr.) Kotl. & Pouzar'
$ws.Range("Q7").Value = 540661.0419420782
$ws.Range("R7").Value = 7247564.172119373

# ---- Row 8 ----
$ws.Range("A8").Value = 111814303
$ws.Range("B8").Value = 90087
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 3298
$ws.Range("F8").Value = 'Trådticka'
$ws.Range("G8").Value = 'Climacocystis borealis'
$ws.Range("H8").Value = '(Fr.) Kotl. & Pouzar'
$ws.Range("J8").Value = ""
$ws.Range("L8").Value = ""
$ws.Range("M8").Value = ""
$ws.Range("Q8").Value = 540600.641023421
$ws.Range("R8").Value = 7247517.393825463
$ws.Range("AF8").Value = ""
$ws.Range("AJ8").Value = ""
$ws.Range("AK8").Value = ""
$ws.Range("AO8").Value = ""

# ---- Row 9 ----
$ws.Range("A9").Value = 111813785
$ws.Range("Q9").Value = 540570.9514120822
$ws.Range("R9").Value = 7247577.960198429

# ---- Row 10 ----
$ws.Range("A10").Value = 111814047
$ws.Range("B10").Value = 90087
$ws.Range("D10").Value = 'LC'
$ws.Range("E10").Value = 3298
$ws.Range("F10").Value = 'Trådticka'
$ws.Range("G10").Value = 'Climacocystis borealis'
$ws.Range("H10").Value = '(Fr.) Kotl. & Pouzar'
$ws.Range("J10").Value = ""
$ws.Range("L10").Value = ""
$ws.Range("M10").Value = ""
$ws.Range("Q10").Value = 540633.6855369165
$ws.Range("R10").Value = 7247516.598344535
$ws.Range("AC10").Value = 'rikligt'
$ws.Range("AF10").Value = ""
$ws.Range("AJ10").Value = ""
$ws.Range("AK10").Value = ""
$ws.Range("AO10").Value = ""

# ---- Row 11 ----
$ws.Range("A11").Value = 111814212
$ws.Range("B11").Value = 89405
$ws.Range("E11").Value = 1202
$ws.Range("F11").Value = 'Ullticka'
$ws.Range("G11").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H11").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q11").Value = 540635.9369002836
$ws.Range("R11").Value = 7247595.565451766

# ---- Row 12 ----
$ws.Range("A12").Value = 111814152
$ws.Range("B12").Value = 89423
$ws.Range("D12").Value = 'NT'
$ws.Range("E12").Value = 5432
$ws.Range("F12").Value = 'Granticka'
$ws.Range("G12").Value = 'Porodaedalea chrysoloma'
$ws.Range("H12").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q12").Value = 540661.0419420782
$ws.Range("R12").Value = 7247564.172119373
$ws.Range("AC12").Value = ""

# ---- Row 13 ----
$ws.Range("A13").Value = 111813745
$ws.Range("B13").Value = 56398
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = 'Tretåig hackspett'
$ws.Range("G13").Value = 'Picoides tridactylus'
$ws.Range("H13").Value = '(Linnaeus, 1758)'
$ws.Range("J13").Value = ""
$ws.Range("L13").Value = ""
$ws.Range("M13").Value = 'färska spår'
$ws.Range("Q13").Value = 540568.950047517
$ws.Range("R13").Value = 7247601.73830481
$ws.Range("AF13").Value = ""
$ws.Range("AJ13").Value = 'gran'
$ws.Range("AK13").Value = 'Picea abies'
$ws.Range("AO13").Value = 'Picea abies'

# ---- Row 14 ----
$ws.Range("A14").Value = 111825245
$ws.Range("B14").Value = 89745
$ws.Range("D14").Value = 'VU'
$ws.Range("E14").Value = 2062
$ws.Range("F14").Value = 'Ulltickeporing'
$ws.Range("G14").Value = 'Skeletocutis brevispora'
$ws.Range("H14").Value = 'Niemelä'
$ws.Range("AC14").Value = 'Färskt exemplar. Kollekt tog och torkades, gulnade.'
$ws.Range("AJ14").Value = 'ullticka'
$ws.Range("AK14").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("AO14").Value = 'Phellinidium ferrugineofuscum'

# ---- Row 15 ----
$ws.Range("A15").Value = 111825340
$ws.Range("B15").Value = 89686
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 658
$ws.Range("F15").Value = 'Rosenticka'
$ws.Range("G15").Value = 'Rhodofomes roseus'
$ws.Range("H15").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("I15").NumberFormat = "@"
$ws.Range("I15").Value = '25'
$ws.Range("AC15").Value = 'Minst 25 rosentickor på denna låga. Ullticka, ulltickeporing, rynkskinn på samma låga'
$ws.Range("AJ15").Value = ""
$ws.Range("AK15").Value = ""
$ws.Range("AO15").Value = ""

# ---- Row 16 ----
$ws.Range("A16").Value = 111825158
$ws.Range("B16").Value = 89405
$ws.Range("E16").Value = 1202
$ws.Range("F16").Value = 'Ullticka'
$ws.Range("G16").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H16").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I16").Value = ""
$ws.Range("AC16").Value = 'Med ulltickeporing'
$ws.Range("AJ16").Value = 'gran'
$ws.Range("AK16").Value = 'Picea abies'
$ws.Range("AO16").Value = 'Picea abies'

# ---- Row 18 ----
$ws.Range("A18").Value = 111939435
$ws.Range("B18").Value = 89423
$ws.Range("E18").Value = 5432
$ws.Range("F18").Value = 'Granticka'
$ws.Range("G18").Value = 'Porodaedalea chrysoloma'
$ws.Range("H18").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("J18").Value = ""
$ws.Range("L18").Value = ""
$ws.Range("M18").Value = ""
$ws.Range("Q18").Value = 540686.0154365924
$ws.Range("R18").Value = 7247572.453681136
$ws.Range("Y18").NumberFormat = "@"
$ws.Range("Y18").Value = '2023-08-14'
$ws.Range("AA18").NumberFormat = "@"
$ws.Range("AA18").Value = '2023-08-14'
$ws.Range("AF18").Value = ""
$ws.Range("AJ18").Value = 'gran'
$ws.Range("AK18").Value = 'Picea abies'
$ws.Range("AL18").Value = 'Gammal gran'
$ws.Range("AO18").Value = 'Picea abies # Gammal gran'

Write-Output "done"